$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 999
$wsExpo.Range("F5").Value = 2696
$wsExpo.Range("F7").Value = 207
$wsExpo.Range("F10").Value = 55
$wsExpo.Range("F11").Value = 2525
$wsExpo.Range("F12").Value = 630

# Sheet "全部类型" (All types) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 999
$wsAll.Range("F6").Value = 2696
$wsAll.Range("F8").Value = 207
$wsAll.Range("F12").Value = 55
$wsAll.Range("F13").Value = 2525
$wsAll.Range("F14").Value = 630
